function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the contact details on the sheet ---
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Normalize the font color used by the Pincode/Phone columns to black ---
$ws.Columns.Item(9).Font.Color = RGB(0, 0, 0)
$ws.Columns.Item(11).Font.Color = RGB(0, 0, 0)

# --- Row heights grew slightly after the edit (19.5 instead of 18.75) ---
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
